$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.816.75"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.559.69"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.481"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.61"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.782.04"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "1.560.72"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.514"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "26.837.14"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.25"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.18"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.15"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.51"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.64"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.99"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0465"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "1.383.09"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.919"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.527"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.995"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("E43").Value = "  +5.01%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.56"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "1.695.75"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.44"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("D50").Value = "0.0₇0978"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0952"
$ws.Range("E51").Value = "  +1.08%  "
